$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.957.23'
$ws.Range("E2").Value = '  -0.61%  '

$ws.Range("D3").Value = '1.638.16'
$ws.Range("E3").Value = '  +0.03%  '

$ws.Range("E4").Value = '  +0.88%  '

$ws.Range("D5").Value = '''214.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.30%  '

$ws.Range("E6").Value = '  +0.56%  '

$ws.Range("E7").Value = '  +0.91%  '

$ws.Range("E8").Value = '  -0.62%  '

$ws.Range("E9").Value = '  +0.66%  '

$ws.Range("D10").Value = '''19.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.88%  '

$ws.Range("E11").Value = '  +0.94%  '

$ws.Range("D12").Value = '1.864.75'
$ws.Range("E12").Value = '  -0.02%  '

$ws.Range("E13").Value = '  -0.12%  '

$ws.Range("D14").Value = '1.638.14'
$ws.Range("E14").Value = '  +0.25%  '

$ws.Range("E15").Value = '  -1.65%  '

$ws.Range("D16").Value = '0.0₃0760'
$ws.Range("E16").Value = '  -0.60%  '

$ws.Range("D17").Value = '''62.56'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.33%  '

$ws.Range("D18").Value = '25.956.93'
$ws.Range("E18").Value = '  -0.55%  '

$ws.Range("E19").Value = '  +0.95%  '

$ws.Range("D20").Value = '''194.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.23%  '

$ws.Range("D21").Value = '''4.37'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.16%  '

$ws.Range("E22").Value = '  -0.71%  '

$ws.Range("E23").Value = '  -1.35%  '

$ws.Range("D24").Value = '''143.99'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.76%  '

$ws.Range("E25").Value = '  +0.32%  '

$ws.Range("E26").Value = '  +1.15%  '

$ws.Range("E27").Value = '  +2.25%  '

$ws.Range("E28").Value = '  -0.83%  '

$ws.Range("E29").Value = '  -0.95%  '

$ws.Range("E30").Value = '  -0.19%  '

$ws.Range("E31").Value = '  +0.68%  '

$ws.Range("E32").Value = '  -1.44%  '

$ws.Range("E33").Value = '  -0.31%  '

$ws.Range("D34").Value = '''1.54'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.87%  '

$ws.Range("E35").Value = '  +2.01%  '

$ws.Range("D36").Value = '''0.905'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.35%  '

$ws.Range("D37").Value = '1.139.28'
$ws.Range("E37").Value = '  -0.14%  '

$ws.Range("D38").Value = '''0.545'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.66%  '

$ws.Range("E39").Value = '  -1.46%  '

$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("D41").Value = '''99.38'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.78%  '

$ws.Range("D42").Value = '''0.800'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.76%  '

$ws.Range("E43").Value = '  -2.72%  '

$ws.Range("D44").Value = '1.774.77'
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("E45").Value = '  +6.90%  '

$ws.Range("D46").Value = '''56.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.07%  '

$ws.Range("E47").Value = '  +2.55%  '

$ws.Range("E48").Value = '  -0.61%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''7.66'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.62%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '''0.415'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.07%  '

$ws.Range("E51").Value = '  -1.24%  '
